# Auto-generated Excel COM-interop script to apply cell value updates
# to the "Phantom_Profits" leve-profit tables across all 8 sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9 (Leve Item ID 5487), hunk 0
$ws.Range("H9").Value = 783
$ws.Range("I9").Value = 783
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 783
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -614
$ws.Range("N9").Value = ""
# Row 11 (Leve Item ID 5533), hunk 1
$ws.Range("H11").Value = 23.444445
$ws.Range("I11").Value = 23.444445
$ws.Range("K11").Value = 23.444445
$ws.Range("M11").Value = 116.555555
# Row 19 (Leve Item ID 7015), hunk 2
$ws.Range("H19").Value = 1328
$ws.Range("I19").Value = 1394.4546
$ws.Range("J19").Value = 1236.625
$ws.Range("K19").Value = 1394.4546
$ws.Range("L19").Value = 1236.625
$ws.Range("M19").Value = -1219.4546
$ws.Range("N19").Value = -1586.625
# Row 86 (Leve Item ID 12603), hunk 3
$ws.Range("H86").Value = 3856.4546
$ws.Range("J86").Value = 3682.6
$ws.Range("L86").Value = 3682.6
$ws.Range("N86").Value = -5928.6
# Row 89 (Leve Item ID 12603), hunk 4
$ws.Range("H89").Value = 3856.4546
$ws.Range("J89").Value = 3682.6
$ws.Range("L89").Value = 18413
$ws.Range("N89").Value = -29645
# Row 92 (Leve Item ID 19901), hunk 5
$ws.Range("H92").Value = 306.91666
$ws.Range("I92").Value = 316.8
$ws.Range("J92").Value = 257.5
$ws.Range("K92").Value = 316.8
$ws.Range("L92").Value = 257.5
$ws.Range("M92").Value = 931.2
$ws.Range("N92").Value = -2753.5
# Row 103 (Leve Item ID 19909), hunk 6
$ws.Range("H103").Value = 936.4286
$ws.Range("J103").Value = 900
$ws.Range("L103").Value = 2700
$ws.Range("N103").Value = -3872
# Row 107 (Leve Item ID 27766), hunk 7
$ws.Range("H107").Value = 902.8
$ws.Range("I107").Value = 703.1111
$ws.Range("K107").Value = 703.1111
$ws.Range("M107").Value = 1216.8889

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713), hunk 8
$ws.Range("H2").Value = 2057.1428
$ws.Range("I2").Value = 1100
$ws.Range("K2").Value = 1100
$ws.Range("M2").Value = -987
# Row 54 (Leve Item ID 2817), hunk 9
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = ""
$ws.Range("N54").Value = ""
# Row 74 (Leve Item ID 44000), hunk 10
$ws.Range("H74").Value = 3944.4
$ws.Range("I74").Value = 2374.6667
$ws.Range("K74").Value = 2374.6667
$ws.Range("M74").Value = -1500.6667
# Row 77 (Leve Item ID 44000), hunk 11
$ws.Range("H77").Value = 3944.4
$ws.Range("I77").Value = 2374.6667
$ws.Range("K77").Value = 11873.3335
$ws.Range("M77").Value = -7505.333500000001
# Row 88 (Leve Item ID 12530), hunk 12
$ws.Range("H88").Value = 953.25
$ws.Range("I88").Value = 1102.1428
$ws.Range("J88").Value = 873.0769
$ws.Range("K88").Value = 1102.1428
$ws.Range("L88").Value = 873.0769
$ws.Range("M88").Value = -696.1428000000001
$ws.Range("N88").Value = -1685.0769
# Row 91 (Leve Item ID 12530), hunk 13
$ws.Range("H91").Value = 953.25
$ws.Range("I91").Value = 1102.1428
$ws.Range("J91").Value = 873.0769
$ws.Range("K91").Value = 1102.1428
$ws.Range("L91").Value = 873.0769
$ws.Range("M91").Value = 301.8571999999999
$ws.Range("N91").Value = -3681.0769
# Row 97 (Leve Item ID 19941), hunk 14
$ws.Range("H97").Value = 906.1786
$ws.Range("I97").Value = 906.1786
$ws.Range("K97").Value = 906.1786
$ws.Range("M97").Value = -410.1786
# Row 110 (Leve Item ID 27708), hunk 15
$ws.Range("H110").Value = 1162.7693
$ws.Range("I110").Value = 828.7273
$ws.Range("K110").Value = 828.7273
$ws.Range("M110").Value = 1216.2727
# Row 116 (Leve Item ID 27713), hunk 16
$ws.Range("H116").Value = 2057.1428
$ws.Range("I116").Value = 1100
$ws.Range("K116").Value = 1100
$ws.Range("M116").Value = 1194
# Row 139 (Leve Item ID 42321), hunk 17
$ws.Range("H139").Value = 79000
$ws.Range("J139").Value = 86250
$ws.Range("L139").Value = 86250
$ws.Range("N139").Value = -96530

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713), hunk 18
$ws.Range("H3").Value = 2057.1428
$ws.Range("I3").Value = 1100
$ws.Range("K3").Value = 1100
$ws.Range("M3").Value = -986
# Row 105 (Leve Item ID 19947), hunk 19
$ws.Range("H105").Value = 2319.6
$ws.Range("I105").Value = 2319.6
$ws.Range("K105").Value = 2319.6
$ws.Range("M105").Value = -572.5999999999999
# Row 107 (Leve Item ID 27706), hunk 20
$ws.Range("H107").Value = 1358.3334
$ws.Range("I107").Value = 1296.4375
$ws.Range("K107").Value = 1296.4375
$ws.Range("M107").Value = 623.5625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023), hunk 21
$ws.Range("H31").Value = 13507.083
$ws.Range("I31").Value = 14565.111
$ws.Range("K31").Value = 14565.111
$ws.Range("M31").Value = -14270.111
# Row 34 (Leve Item ID 44023), hunk 22
$ws.Range("H34").Value = 13507.083
$ws.Range("I34").Value = 14565.111
$ws.Range("K34").Value = 14565.111
$ws.Range("M34").Value = -14363.111
# Row 50 (Leve Item ID 1862), hunk 23
$ws.Range("H50").Value = 29610.4
$ws.Range("I50").Value = 19350.666
$ws.Range("J50").Value = 45000
$ws.Range("K50").Value = 19350.666
$ws.Range("L50").Value = 45000
$ws.Range("M50").Value = -18725.666
$ws.Range("N50").Value = -46250
# Row 119 (Leve Item ID 26276), hunk 24
$ws.Range("H119").Value = 39999.3
$ws.Range("J119").Value = 39999.3
$ws.Range("L119").Value = 39999.3
$ws.Range("N119").Value = -49675.3
# Row 121 (Leve Item ID 27227), hunk 25
$ws.Range("H121").Value = 31331.666
$ws.Range("J121").Value = 31331.666
$ws.Range("L121").Value = 31331.666
$ws.Range("N121").Value = -33951.666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 43 (Leve Item ID 4669), hunk 26
$ws.Range("H43").Value = 19899.25
$ws.Range("J43").Value = 19899.25
$ws.Range("L43").Value = 59697.75
$ws.Range("N43").Value = -59925.75
# Row 56 (Leve Item ID 10146), hunk 27
$ws.Range("H56").Value = 9603.733
$ws.Range("I56").Value = 9603.733
$ws.Range("K56").Value = 9603.733
$ws.Range("M56").Value = -9073.733
# Row 60 (Leve Item ID 4750), hunk 28
$ws.Range("H60").Value = 1418.3572
$ws.Range("J60").Value = 2108.75
$ws.Range("L60").Value = 6326.25
$ws.Range("N60").Value = -6828.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102 (Leve Item ID 36169), hunk 29
$ws.Range("H102").Value = 979.4074000000001
$ws.Range("I102").Value = 995.125
$ws.Range("J102").Value = 853.6667
$ws.Range("K102").Value = 995.125
$ws.Range("L102").Value = 853.6667
$ws.Range("M102").Value = 626.875
$ws.Range("N102").Value = -4097.6667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 82 (Leve Item ID 12565), hunk 30
$ws.Range("H82").Value = 1114.5
$ws.Range("I82").Value = 433.5
$ws.Range("K82").Value = 433.5
$ws.Range("M82").Value = -72.5
# Row 85 (Leve Item ID 12565), hunk 31
$ws.Range("H85").Value = 1114.5
$ws.Range("I85").Value = 433.5
$ws.Range("K85").Value = 433.5
$ws.Range("M85").Value = 814.5
# Row 93 (Leve Item ID 19993), hunk 32
$ws.Range("H93").Value = 1017.8571
$ws.Range("I93").Value = 965.2
$ws.Range("K93").Value = 965.2
$ws.Range("M93").Value = 282.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 3 (Leve Item ID 3309), hunk 33
$ws.Range("H3").Value = 6254947
$ws.Range("J3").Value = 6596.3335
$ws.Range("L3").Value = 6596.3335
$ws.Range("N3").Value = -6824.3335
# Row 53 (Leve Item ID 3172), hunk 34
$ws.Range("H53").Value = 10000
$ws.Range("I53").Value = 10000
$ws.Range("K53").Value = 10000
$ws.Range("M53").Value = -9393
# Row 62 (Leve Item ID 12589), hunk 35
$ws.Range("H62").Value = 6159.8
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 6699.75
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 6699.75
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -7947.75
# Row 65 (Leve Item ID 12589), hunk 36
$ws.Range("H65").Value = 6159.8
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 6699.75
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 33498.75
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -39738.75
# Row 141 (Leve Item ID 42505), hunk 37
$ws.Range("H141").Value = 180325
$ws.Range("I141").Value = 210650
$ws.Range("J141").Value = 150000
$ws.Range("K141").Value = 210650
$ws.Range("L141").Value = 150000
$ws.Range("M141").Value = -205470
$ws.Range("N141").Value = -160360

